# Cambiar nombre de la etiqueta "FORMATO INTERMEDIO" de la entidad
# Servicio por "FORMATO ENTRADA/SALIDA" (hoja "Cartera", columna X,
# fila de cabecera 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cartera")

$ws.Range("X1").Value = "FORMATO ENTRADA/SALIDA"

# El texto nuevo es mas largo que el anterior: ajustamos el ancho de la
# columna X para que se ajuste al contenido (equivalente a un "autofit").
$ws.Columns.Item(24).ColumnWidth = 25.86

# Desplazar la vista de la hoja hacia la derecha (la nueva columna de
# interes queda mas a la derecha) y actualizar la celda seleccionada.
$excel.ActiveWindow.ScrollColumn = 18
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("W3").Select()

# Ajustar el ancho de la ventana del libro.
$excel.ActiveWindow.Width = 19440
